# Handles float input without breaking stuff
# Fills in the "Student Ans" column with the student's actual answers,
# recomputes the right/wrong/not-attempted summary block, and drops the
# two unused duplicate "Student Ans / Correct Ans" blocks (columns D:E
# beyond the header row, and columns G:H entirely).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

function Set-StudentAnswer($row, $value, $styleName) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $value
    $cell.Style = $styleName
    $cell.HorizontalAlignment = $xlCenter
}

# ---- Summary block (rows 10-12) ----
# Give the row-label cells in column A the same "mtitleStyle" formatting
# used by the row above (row 9).
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A10").HorizontalAlignment = $xlCenter
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A11").HorizontalAlignment = $xlCenter
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("A12").HorizontalAlignment = $xlCenter

# Row 10 "No." - counts of right / wrong / not attempted / max
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 28

# Row 11 "Marking" - marks per right / wrong / not attempted
# (wrong-answer penalty now stored as a real number instead of text)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 "Total"
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -9
$ws.Range("E12").Value = "55/112"

# ---- Per-question "Student Ans" column (A16:A40) ----
# Column B ("Correct Ans") already holds the answer key; fill in column A
# with what the student actually answered, colouring it green
# (correctStyle) when it matches column B and red (incorrectStyle) when
# it doesn't. Rows 24, 29 and 34 were left unanswered and already carry
# the default blank styling, so they're untouched.
Set-StudentAnswer 16 "Option A" "correctStyle"
Set-StudentAnswer 17 "Option D" "correctStyle"
Set-StudentAnswer 18 "Option D" "incorrectStyle"
Set-StudentAnswer 19 "Option C" "correctStyle"
Set-StudentAnswer 20 "Option B" "correctStyle"
Set-StudentAnswer 21 "Option B" "incorrectStyle"
Set-StudentAnswer 22 "Option D" "correctStyle"
Set-StudentAnswer 23 "Option D" "correctStyle"
Set-StudentAnswer 25 "Option A" "correctStyle"
Set-StudentAnswer 26 "Option C" "correctStyle"
Set-StudentAnswer 27 "Option B" "incorrectStyle"
Set-StudentAnswer 28 "Option D" "correctStyle"
Set-StudentAnswer 30 "Option B" "correctStyle"
Set-StudentAnswer 31 "Option B" "incorrectStyle"
Set-StudentAnswer 32 "Option C" "correctStyle"
Set-StudentAnswer 33 "Option A" "incorrectStyle"
Set-StudentAnswer 35 "Option D" "correctStyle"
Set-StudentAnswer 36 "Option C" "incorrectStyle"
Set-StudentAnswer 37 "Option A" "correctStyle"
Set-StudentAnswer 38 "Option A" "correctStyle"
Set-StudentAnswer 39 "Option B" "incorrectStyle"
Set-StudentAnswer 40 "Option D" "correctStyle"

# ---- Drop the unused duplicate answer blocks ----
# Columns D:E carried a second (unused) "Student Ans / Correct Ans" block
# below the row-15 header; columns G:H carried a third one including its
# own header. Both are fully cleared (formatting included) so the sheet's
# used range shrinks back down to A:E.
$ws.Range("D16:E40").Clear()
$ws.Range("G15:H40").Clear()

Write-Host "edit complete"
